$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59: clear C59 (remove the "NA" text, leave it blank)
$ws.Range("C59").Value = $null

# Row 60
$ws.Range("A60").Value = "'2025-04-28"
$ws.Range("B60").Value = "bonnes pratiques"
$ws.Range("C60").Value = 45
$ws.Range("D60").Value = 4

# Row 61
$ws.Range("A61").Value = "'2025-04-28"
$ws.Range("B61").Value = "bonnes pratiques"
$ws.Range("C61").Value = 46
$ws.Range("D61").Value = 6

# Row 62
$ws.Range("A62").Value = "'2025-04-28"
$ws.Range("B62").Value = "bonnes pratiques"
$ws.Range("C62").Value = 49
$ws.Range("D62").Value = 1
